$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11: employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(2,  21973, "Hadassa Vieira",      "Juridico",       "Problemas pessoais", 5, 45104, 3106.61),
    @(3,  45472, "Dom Nascimento",      "P&D",            "Problemas pessoais", 5, 45084, 3759.98),
    @(4,  91486, "Rael Guerra",         "Vendas",         "Outros",             2, 45103, 3941.37),
    @(5,  26843, "Enzo Dias",           "Financeiro",     "Consulta medica",    4, 45097, 4454.95),
    @(6,  45551, "Sra. Cecília Novaes", "Marketing",      "Consulta medica",    2, 45084, 3964.51),
    @(7,  69038, "Samuel Barbosa",      "Recursos Humanos","Problemas pessoais",5, 45081, 6895.21),
    @(8,  34882, "Yuri da Paz",         "Operacoes",      "Outros",             5, 45099, 6186.1),
    @(9,  98795, "Joana Fonseca",       "TI",             "Outros",             2, 45104, 5313.39),
    @(10, 40032, "Rafaela Pinto",       "Marketing",      "Doenca",             5, 45091, 8694.74),
    @(11, 14655, "Aurora Marques",      "Juridico",       "Consulta medica",    7, 45100, 6353.67)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
